$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------------
# Sheet1 header rework.
#   old: A1=EDSNumber, B1=Reason,               C1=AttemptDate
#   new: A1=EDSNumber, B1=FirstAttemptReason,    C1=FirstAttemptDate,
#                       D1=SecondAttemptReason,   E1=SecondAttemptDate
# Write right-to-left (E,D,C,B) so the shared-string table is appended in
# the same order the original authoring session produced it.
# ---------------------------------------------------------------------------
$ws1.Range("E1").Value = "SecondAttemptDate"
$ws1.Range("D1").Value = "SecondAttemptReason"
$ws1.Range("C1").Value = "FirstAttemptDate"
$ws1.Range("B1").Value = "FirstAttemptReason"

# D1 is a plain centered header cell, like B1.
$ws1.Range("D1").HorizontalAlignment = -4108

# E1 is a centered date cell, using the same custom date format as C1.
$ws1.Range("E1").NumberFormat = "[$-409]d\-mmm\-yyyy;@"
$ws1.Range("E1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Data validation
# ---------------------------------------------------------------------------
# E1:E1048576 gets the same "type DD/MM/YYYY" input-message validation as C1.
$ws1.Range("E1:E1048576").Validation.Add(0)
$ws1.Range("E1:E1048576").Validation.InputTitle = "Perhatikan Format Penulisan"
$ws1.Range("E1:E1048576").Validation.InputMessage = "Ketik: DD/MM/YYYY_x000a_Contoh: 16/12/2016"

# D2:D1048576 gets a new list validation sourced from Sheet2!$A$1:$A$12
# (a bounded version of the Sheet2!$A:$A list already used for column B).
$ws1.Range("D2:D1048576").Validation.Add(3, 1, 1, "Sheet2!`$A`$1:`$A`$12")

# ---------------------------------------------------------------------------
# Column widths - mirror Excel's "autofit" result for the new header text
# (closest attainable quantized widths for this engine's width model).
# ---------------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 10.0
$ws1.Columns.Item(2).ColumnWidth = 16.5
$ws1.Columns.Item(3).ColumnWidth = 14.5
$ws1.Columns.Item(4).ColumnWidth = 18.833333333333332
$ws1.Columns.Item(5).ColumnWidth = 16.833333333333332

# ---------------------------------------------------------------------------
# Selections - Sheet2 lands on A7, Sheet1 (the active tab) lands on B7.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A7").Select()
$ws1.Activate()
$ws1.Range("B7").Select()
